$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("F:F").Insert()
$ws.Range("F1").Value = "Plate_Barcode"
